$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
# Values that look like plain numbers (single '.' decimal point) are written with a
# leading apostrophe so Excel keeps them as literal text, matching the source data's
# inline-string cells (prices such as "1.000" must not collapse to the number 1).

$ws.Range("D2").Value = "29.871.59"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.871.12"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'0.7325"
$ws.Range("E5").Value = "  -5.77%  "
$ws.Range("D6").Value = "'241.92"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.3152"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").Value = "'24.66"
$ws.Range("E9").Value = "  -4.89%  "
$ws.Range("D10").Value = "'0.07089"
$ws.Range("D11").Value = "'0.08422"
$ws.Range("E11").Value = "  -8.64%  "
$ws.Range("D12").Value = "'0.7512"
$ws.Range("E12").Value = "  -2.90%  "
$ws.Range("D13").Value = "'5.394"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").Value = "1.864.68"
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("D15").Value = "'92.52"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "29.883.60"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "'6.039"
$ws.Range("E17").Value = "  -3.17%  "
$ws.Range("E18").Value = "  -3.00%  "
$ws.Range("D19").Value = "'242.71"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("D20").Value = "'0.000007818"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "2.121.21"
$ws.Range("E22").Value = "  -3.57%  "
$ws.Range("D23").Value = "'7.922"
$ws.Range("E23").Value = "  -3.26%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -2.42%  "
$ws.Range("D26").Value = "'9.305"
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("D27").Value = "'163.83"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").Value = "'18.57"
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").Value = "'2.015"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("D30").Value = "'1.471"
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("D31").Value = "'4.583"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D32").Value = "'1.530"
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("D33").Value = "'4.275"
$ws.Range("E33").Value = "  +3.69%  "
$ws.Range("D34").Value = "'0.05330"
$ws.Range("E34").Value = "  -3.19%  "
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").Value = "'0.7529"
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").Value = "'1.000"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "'2.702"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").Value = "'0.01951"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "'2.749"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").Value = "'0.4468"
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("D42").Value = "1.107.19"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Value = "'6.073"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "'72.18"
$ws.Range("E44").Value = "  -2.92%  "
$ws.Range("D45").Value = "'0.8648"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "'102.81"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("D51").Value = "2.019.78"
$ws.Range("E51").Value = "  -3.22%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.840"
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "'3.064"
$ws.Range("E50").Value = "  +2.03%  "
